$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 4168903.8
$ws.Range("I40").Value = 31250900
$ws.Range("J40").Value = 2442.923
$ws.Range("K40").Value = 31250900
$ws.Range("L40").Value = 2442.923
$ws.Range("M40").Value = -31250725
$ws.Range("N40").Value = -2792.923

$ws.Range("H129").Value = 17242392
$ws.Range("I129").Value = 111112180
$ws.Range("J129").Value = 1003.4286
$ws.Range("K129").Value = 333336540
$ws.Range("L129").Value = 3010.2858
$ws.Range("M129").Value = -333331540
$ws.Range("N129").Value = -13010.2858

$ws.Range("H131").Value = 3425
$ws.Range("I131").Value = 817
$ws.Range("K131").Value = 2451
$ws.Range("M131").Value = 2589

$ws.Range("H137").Value = 1452.2333
$ws.Range("I137").Value = 1202.3043
$ws.Range("J137").Value = 2273.4285
$ws.Range("K137").Value = 3606.9129
$ws.Range("L137").Value = 6820.2855
$ws.Range("M137").Value = -1056.9129
$ws.Range("N137").Value = -11920.2855

$ws.Range("H138").Value = 3756.5
$ws.Range("I138").Value = 2386.6086
$ws.Range("J138").Value = 4241.231
$ws.Range("K138").Value = 7159.825800000001
$ws.Range("L138").Value = 12723.693
$ws.Range("M138").Value = -2019.825800000001
$ws.Range("N138").Value = -23003.693

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 30860.412
$ws.Range("I45").Value = 35835.93
$ws.Range("J45").Value = 2002.4
$ws.Range("K45").Value = 35835.93
$ws.Range("L45").Value = 2002.4
$ws.Range("M45").Value = -35458.93
$ws.Range("N45").Value = -2756.4

$ws.Range("H55").Value = 47780
$ws.Range("J55").Value = 47780
$ws.Range("L55").Value = 47780
$ws.Range("N55").Value = -48410

$ws.Range("H61").Value = 1772.9841
$ws.Range("I61").Value = 1702.3513
$ws.Range("J61").Value = 1873.5
$ws.Range("K61").Value = 1702.3513
$ws.Range("L61").Value = 1873.5
$ws.Range("M61").Value = -1490.3513
$ws.Range("N61").Value = -2297.5

$ws.Range("H74").Value = 1796.0625
$ws.Range("I74").Value = 2017.7878
$ws.Range("J74").Value = 1308.2667
$ws.Range("K74").Value = 2017.7878
$ws.Range("L74").Value = 1308.2667
$ws.Range("M74").Value = -1143.7878
$ws.Range("N74").Value = -3056.2667

$ws.Range("H77").Value = 1796.0625
$ws.Range("I77").Value = 2017.7878
$ws.Range("J77").Value = 1308.2667
$ws.Range("K77").Value = 10088.939
$ws.Range("L77").Value = 6541.3335
$ws.Range("M77").Value = -5720.939
$ws.Range("N77").Value = -15277.3335

$ws.Range("H136").Value = 1772.9841
$ws.Range("I136").Value = 1702.3513
$ws.Range("J136").Value = 1873.5
$ws.Range("K136").Value = 5107.0539
$ws.Range("L136").Value = 5620.5
$ws.Range("M136").Value = -2557.0539
$ws.Range("N136").Value = -10720.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2525.4666
$ws.Range("I20").Value = 2570.611
$ws.Range("J20").Value = 2457.75
$ws.Range("K20").Value = 2570.611
$ws.Range("L20").Value = 2457.75
$ws.Range("M20").Value = -2323.611
$ws.Range("N20").Value = -2951.75

$ws.Range("H63").Value = 17219
$ws.Range("J63").Value = 17219
$ws.Range("L63").Value = 17219
$ws.Range("N63").Value = -18591

$ws.Range("H66").Value = 17219
$ws.Range("J66").Value = 17219
$ws.Range("L66").Value = 51657
$ws.Range("N66").Value = -58521

$ws.Range("H134").Value = 3182801
$ws.Range("I134").Value = 6115.64
$ws.Range("J134").Value = 11124514
$ws.Range("K134").Value = 18346.92
$ws.Range("L134").Value = 33373542
$ws.Range("M134").Value = -15811.92
$ws.Range("N134").Value = -33378612

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 266.83334
$ws.Range("I22").Value = 250.25
$ws.Range("J22").Value = 300
$ws.Range("K22").Value = 250.25
$ws.Range("L22").Value = 300
$ws.Range("M22").Value = 99.75
$ws.Range("N22").Value = -1000

$ws.Range("H31").Value = 2326.9883
$ws.Range("I31").Value = 1561.525
$ws.Range("J31").Value = 2992.6086
$ws.Range("K31").Value = 1561.525
$ws.Range("L31").Value = 2992.6086
$ws.Range("M31").Value = -1266.525
$ws.Range("N31").Value = -3582.6086

$ws.Range("H34").Value = 2326.9883
$ws.Range("I34").Value = 1561.525
$ws.Range("J34").Value = 2992.6086
$ws.Range("K34").Value = 1561.525
$ws.Range("L34").Value = 2992.6086
$ws.Range("M34").Value = -1359.525
$ws.Range("N34").Value = -3396.6086

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 71430184
$ws.Range("I16").Value = 6494590.5
$ws.Range("K16").Value = 6494590.5
$ws.Range("M16").Value = -6494420.5

$ws.Range("H22").Value = 1267342.9
$ws.Range("I22").Value = 4219609
$ws.Range("J22").Value = 2086
$ws.Range("K22").Value = 4219609
$ws.Range("L22").Value = 2086
$ws.Range("M22").Value = -4219314
$ws.Range("N22").Value = -2676

$ws.Range("H27").Value = 1267342.9
$ws.Range("I27").Value = 4219609
$ws.Range("J27").Value = 2086
$ws.Range("K27").Value = 4219609
$ws.Range("L27").Value = 2086
$ws.Range("M27").Value = -4219502
$ws.Range("N27").Value = -2300

$ws.Range("H43").Value = 4000
$ws.Range("J43").Value = 4000
$ws.Range("L43").Value = 4000
$ws.Range("N43").Value = -4386

$ws.Range("H46").Value = 5209009.5
$ws.Range("I46").Value = 13889454
$ws.Range("J46").Value = 742.4
$ws.Range("K46").Value = 13889454
$ws.Range("L46").Value = 742.4
$ws.Range("M46").Value = -13889266
$ws.Range("N46").Value = -1118.4

$ws.Range("H55").Value = 166683470
$ws.Range("I55").Value = 33460
$ws.Range("J55").Value = 333333470
$ws.Range("K55").Value = 33460
$ws.Range("L55").Value = 333333470
$ws.Range("M55").Value = -33287
$ws.Range("N55").Value = -333333816

$ws.Range("H61").Value = 1645.7084
$ws.Range("I61").Value = 1379.3158
$ws.Range("J61").Value = 2658
$ws.Range("K61").Value = 1379.3158
$ws.Range("L61").Value = 2658
$ws.Range("M61").Value = -1177.3158
$ws.Range("N61").Value = -3062

$ws.Range("H113").Value = 1645.7084
$ws.Range("I113").Value = 1379.3158
$ws.Range("J113").Value = 2658
$ws.Range("K113").Value = 1379.3158
$ws.Range("L113").Value = 2658
$ws.Range("M113").Value = 790.6841999999999
$ws.Range("N113").Value = -6998

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2181.6667
$ws.Range("I122").Value = 1492.6842
$ws.Range("J122").Value = 4799.8
$ws.Range("K122").Value = 4478.0526
$ws.Range("L122").Value = 14399.4
$ws.Range("M122").Value = -2028.0526
$ws.Range("N122").Value = -19299.4

$ws.Range("H132").Value = 14943943
$ws.Range("I132").Value = 21299892
$ws.Range("K132").Value = 63899676
$ws.Range("M132").Value = -63897146
